# Auto-generated Excel COM-interop script to apply the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (e.g. "1.00", "27.18") are stored as literal text,
# matching the source data which keeps these as plain strings rather than numbers.
$textCells = @("D4","D5","D6","D7","D10","D13","D14","D15","D18","D20","D21","D23","D24","D25","D26","D27","D30","D33","D35","D37","D38","D39","D40","D41","D42","D44","D45","D48","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '64.951.93'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '2.949.96'
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '569.22'
$ws.Range("E5").Value = '  -2.51%  '
$ws.Range("D6").Value = '159.03'
$ws.Range("E6").Value = '  +2.97%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("E8").Value = '  +0.69%  '
$ws.Range("D9").Value = '2.945.70'
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").Value = '6.67'
$ws.Range("E10").Value = '  -5.11%  '
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  +2.09%  '
$ws.Range("D14").Value = '34.12'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '0.125'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '65.036.20'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = '3.439.38'
$ws.Range("E17").Value = '  -1.42%  '
$ws.Range("D18").Value = '6.94'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").Value = '2.949.67'
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").Value = '446.00'
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("D21").Value = '13.85'
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '7.25'
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").Value = '82.48'
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("D25").Value = '2.21'
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("D26").Value = '12.02'
$ws.Range("E26").Value = '  -3.73%  '
$ws.Range("D27").Value = '10.10'
$ws.Range("E27").Value = '  -6.21%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  +2.63%  '
$ws.Range("D30").Value = '2.37'
$ws.Range("E30").Value = '  -2.68%  '
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("E32").Value = '  -2.75%  '
$ws.Range("D33").Value = '27.18'
$ws.Range("E33").Value = '  +1.22%  '
$ws.Range("E34").Value = '  -1.60%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").Value = '5.68'
$ws.Range("E37").Value = '  -1.97%  '
$ws.Range("D38").Value = '48.92'
$ws.Range("E38").Value = '  -0.47%  '
$ws.Range("D39").Value = '1.97'
$ws.Range("E39").Value = '  -6.28%  '
$ws.Range("D40").Value = '43.96'
$ws.Range("E40").Value = '  -4.89%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").Value = '0.299'
$ws.Range("E41").Value = '  -1.39%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '2.83'
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("E43").Value = '  -1.67%  '
$ws.Range("D44").Value = '8.40'
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '384.81'
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").Value = '2.719.28'
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("D48").Value = '132.24'
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").Value = '2.15'
$ws.Range("E50").Value = '  +4.59%  '
$ws.Range("E51").Value = '  +0.63%  '
